$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Is" (column B) is now a constant 10 A for every row (was 9.5 A)
$ws.Range("B7:B11").Value = 10

# "Ilmax" (column A) rows re-ordered into an ascending 12/15/20/25/30 sequence
$ws.Range("A7").Value = 12
$ws.Range("A8").Value = 15
$ws.Range("A9").Value = 20
$ws.Range("A10").Value = 25
$ws.Range("A11").Value = 30

# Recompute the derived columns (C:H) for the reshuffled rows 8:11 so the
# dependent formulas/values stay in sync with the new A/B inputs.
$ws.Range("C8:H11").ClearContents()
$ws.Range("C8:C11").Formula = "=B8/A8"
$ws.Range("D8:D11").Formula = "=(1-C8)*B8/(2*`$B`$1*(`$B`$3-`$B`$5))"
$ws.Range("E8:E11").Formula = "=2*`$B`$1*D8/(1-C8)"
$ws.Range("F8:F11").Formula = "=(1-C8)*(2*B8/`$B`$3)/(2*`$B`$1)"
$ws.Range("G8:G11").Formula = "=C8/`$B`$1"
$ws.Range("H8:H11").Formula = "=`$B`$2-G8"

# Window/selection state as left by the author: zoomed to 140%, cursor on D12
$excel.ActiveWindow.Zoom = 140
$ws.Range("D12").Select()
